# Update countries & provincias Spain
# The source COVID dashboard data was refreshed: the footer timestamp moved
# from 01:18 to 02:35, several countries' metrics were updated, and because
# the sheet is rank-ordered by case count, a handful of rows swapped which
# country they show (their whole row - name + figures - moved together).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Footer timestamp cell (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Julio de 2020 a las 02:35"

$changes = @(
    @{ Row=4;   Country="Estados Unidos";               NameChanged=$false; Vals=@(4248203,77885,2027615,2072110,0,1129,148478) },
    @{ Row=21;  Country="Alemania";                      NameChanged=$false; Vals=@(205960,818,190400,6359,0,14,9201) },
    @{ Row=24;  Country="Canada";                        NameChanged=$false; Vals=@(113206,534,98873,5452,0,7,8881) },
    @{ Row=43;  Country="Panama";                        NameChanged=$false; Vals=@(57993,1176,32704,24039,0,41,1250) },
    @{ Row=57;  Country="Ghana";                         NameChanged=$false; Vals=@(31057,1385,27801,3095,0,8,161) },
    @{ Row=102; Country="Paraguay";                      NameChanged=$true;  Vals=@(4224,111,2596,1590,0,2,38) },
    @{ Row=103; Country="Grecia";                        NameChanged=$true;  Vals=@(4135,25,1374,2560,0,0,201) },
    @{ Row=105; Country="Malaui";                        NameChanged=$true;  Vals=@(3453,151,1373,1993,0,11,87) },
    @{ Row=106; Country="Nicaragua";                     NameChanged=$true;  Vals=@(3439,0,2492,839,0,0,108) },
    @{ Row=107; Country="Libano";                        NameChanged=$true;  Vals=@(3407,147,1666,1695,0,3,46) },
    @{ Row=140; Country="Uruguay";                       NameChanged=$true;  Vals=@(1166,25,946,186,0,0,34) },
    @{ Row=141; Country="Jordania";                      NameChanged=$true;  Vals=@(1146,15,1035,100,0,0,11) },
    @{ Row=164; Country="Guyana";                        NameChanged=$false; Vals=@(352,1,178,154,0,1,20) },
    @{ Row=170; Country="Bahamas";                       NameChanged=$true;  Vals=@(316,42,91,214,0,0,11) },
    @{ Row=171; Country="Mongolia";                      NameChanged=$true;  Vals=@(288,1,217,71,0,0,0) },
    @{ Row=175; Country="Guadalupe";                     NameChanged=$true;  Vals=@(203,8,176,13,0,0,14) },
    @{ Row=176; Country="Islas Caimanes";                NameChanged=$true;  Vals=@(203,0,202,0,0,0,1) },
    @{ Row=177; Country="Camboya";                       NameChanged=$true;  Vals=@(202,4,142,60,0,0,0) },
    @{ Row=193; Country="San Vicente y las Granadinas";  NameChanged=$false; Vals=@(52,0,39,13,0,0,0) },
    @{ Row=210; Country="Islas Malvinas";                NameChanged=$true;  Vals=@(13,0,13,0,0,0,0) },
    @{ Row=211; Country="Groenlandia";                   NameChanged=$true;  Vals=@(13,0,13,0,0,0,0) }
)

foreach ($chg in $changes) {
    $r = $chg.Row
    if ($chg.NameChanged) {
        $ws.Cells.Item($r, 1).Value = $chg.Country
    }
    $col = 2
    foreach ($v in $chg.Vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}
